# Update column G ("K" = strikeouts) values on the active sheet.
# Replaces the previous "Strike#" derived values with the regenerated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 5
    4  = 5
    5  = 1
    6  = 4
    7  = 0
    8  = 1
    9  = 5
    10 = 4
    11 = 3
    12 = 1
    13 = 1
    14 = 3
    15 = 0
    16 = 3
    17 = 3
    18 = 4
    19 = 3
    20 = 6
    21 = 5
    22 = 4
    23 = 5
    24 = 2
    25 = 3
    26 = 3
    27 = 2
    28 = 6
    29 = 3
    30 = 5
    31 = 4
    32 = 5
    33 = 0
    34 = 7
    35 = 3
    36 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
